$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 17216.648
$ws.Range("I28").Value = 27269.39
$ws.Range("J28").Value = 701.4286
$ws.Range("K28").Value = 27269.39
$ws.Range("L28").Value = 701.4286
$ws.Range("M28").Value = -26784.39
$ws.Range("N28").Value = -1671.4286

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 3829.3635
$ws.Range("I69").Value = 3760
$ws.Range("J69").Value = 3849.7646
$ws.Range("K69").Value = 11280
$ws.Range("L69").Value = 11549.2938
$ws.Range("M69").Value = -10406
$ws.Range("N69").Value = -13297.2938

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H72").Value = 3829.3635
$ws.Range("I72").Value = 3760
$ws.Range("J72").Value = 3849.7646
$ws.Range("K72").Value = 33840
$ws.Range("L72").Value = 34647.8814
$ws.Range("M72").Value = -29472
$ws.Range("N72").Value = -43383.8814

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 8546.954
$ws.Range("I45").Value = 11652.2
$ws.Range("J45").Value = 1892.8572
$ws.Range("K45").Value = 11652.2
$ws.Range("L45").Value = 1892.8572
$ws.Range("M45").Value = -11275.2
$ws.Range("N45").Value = -2646.8572

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1611.3334
$ws.Range("I102").Value = 1417.375
$ws.Range("J102").Value = 1999.25
$ws.Range("K102").Value = 1417.375
$ws.Range("L102").Value = 1999.25
$ws.Range("M102").Value = 204.625
$ws.Range("N102").Value = -5243.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H87").Value = 37088.5
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 37088.5
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 37088.5
$ws.Range("N87").Value = -39584.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H90").Value = 37088.5
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 37088.5
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 111265.5
$ws.Range("N90").Value = -123745.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1064.025
$ws.Range("I134").Value = 904.5862
$ws.Range("J134").Value = 1484.3636
$ws.Range("K134").Value = 2713.7586
$ws.Range("L134").Value = 4453.0908
$ws.Range("M134").Value = -178.7586000000001
$ws.Range("N134").Value = -9523.0908

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 18786.754
$ws.Range("I31").Value = 1102.3611
$ws.Range("J31").Value = 40739.793
$ws.Range("K31").Value = 1102.3611
$ws.Range("L31").Value = 40739.793
$ws.Range("M31").Value = -807.3611000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 18786.754
$ws.Range("I34").Value = 1102.3611
$ws.Range("J34").Value = 40739.793
$ws.Range("K34").Value = 1102.3611
$ws.Range("L34").Value = 40739.793
$ws.Range("M34").Value = -900.3611000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 174.33333
$ws.Range("I33").Value = 56.2
$ws.Range("J33").Value = 233.4
$ws.Range("K33").Value = 337.2
$ws.Range("L33").Value = 1400.4
$ws.Range("M33").Value = -54.20000000000005
$ws.Range("N33").Value = -1966.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 851.4815
$ws.Range("I122").Value = 422.6
$ws.Range("J122").Value = 2076.8572
$ws.Range("K122").Value = 3803.4
$ws.Range("L122").Value = 18691.7148
$ws.Range("M122").Value = -1353.4
$ws.Range("N122").Value = -23591.7148

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 779
$ws.Range("I131").Value = 509.77777
$ws.Range("J131").Value = 985.21277
$ws.Range("K131").Value = 1529.33331
$ws.Range("L131").Value = 2955.63831
$ws.Range("M131").Value = 3510.66669
$ws.Range("N131").Value = -13035.63831

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2369.3572
$ws.Range("I97").Value = 2347.5833
$ws.Range("J97").Value = 2500
$ws.Range("K97").Value = 2347.5833
$ws.Range("L97").Value = 2500
$ws.Range("M97").Value = -1851.5833
$ws.Range("N97").Value = -3492

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1700
$ws.Range("I68").Value = 1800
$ws.Range("J68").Value = 1650
$ws.Range("K68").Value = 1800
$ws.Range("L68").Value = 1650
$ws.Range("M68").Value = -1051
$ws.Range("N68").Value = -3148

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 1700
$ws.Range("I71").Value = 1800
$ws.Range("J71").Value = 1650
$ws.Range("K71").Value = 9000
$ws.Range("L71").Value = 8250
$ws.Range("M71").Value = -5256
$ws.Range("N71").Value = -15738

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H75").Value = 48000
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 48000
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 48000
$ws.Range("N75").Value = -49872
$ws.Range("M75").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H78").Value = 48000
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 48000
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 144000
$ws.Range("N78").Value = -153360
$ws.Range("M78").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3001.5
$ws.Range("I62").Value = 3001.2
$ws.Range("J62").Value = 3003
$ws.Range("K62").Value = 3001.2
$ws.Range("L62").Value = 3003
$ws.Range("M62").Value = -2377.2
$ws.Range("N62").Value = -4251

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 3001.5
$ws.Range("I65").Value = 3001.2
$ws.Range("J65").Value = 3003
$ws.Range("K65").Value = 15006
$ws.Range("L65").Value = 15015
$ws.Range("M65").Value = -11886
$ws.Range("N65").Value = -21255

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H68").Value = 39000
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 39000
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 39000
$ws.Range("N68").Value = -40622

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 48000
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 48000
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 48000
$ws.Range("N69").Value = -49498

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H71").Value = 39000
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 39000
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 117000
$ws.Range("N71").Value = -125112

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H72").Value = 48000
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 48000
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 144000
$ws.Range("N72").Value = -151488

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1144.4445
$ws.Range("I81").Value = 1053.3334
$ws.Range("J81").Value = 1326.6666
$ws.Range("K81").Value = 2106.6668
$ws.Range("L81").Value = 2653.3332
$ws.Range("M81").Value = -1045.6668
$ws.Range("N81").Value = -4775.3332

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 1144.4445
$ws.Range("I84").Value = 1053.3334
$ws.Range("J84").Value = 1326.6666
$ws.Range("K84").Value = 10533.334
$ws.Range("L84").Value = 13266.666
$ws.Range("M84").Value = -5229.333999999999
$ws.Range("N84").Value = -23874.666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H88").Value = 10171
$ws.Range("I88").Value = 10171
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 10171
$ws.Range("L88").Value = 0
$ws.Range("M88").Value = -9765

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H91").Value = 10171
$ws.Range("I91").Value = 10171
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 10171
$ws.Range("L91").Value = 0
$ws.Range("M91").Value = -8767
